# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Agosto de 2020 a las 16:09"

# --- Updated country statistics (B:H columns) ---

# row 4: Estados Unidos
$ws.Range("B4").Value = 5532270
$ws.Range("C4").Value = 2481
$ws.Range("D4").Value = 2904245
$ws.Range("E4").Value = 2455395

# row 6: India
$ws.Range("B6").Value = 2618877
$ws.Range("C6").Value = 29669
$ws.Range("D6").Value = 1891084
$ws.Range("E6").Value = 677329
$ws.Range("G6").Value = 380
$ws.Range("H6").Value = 50464

# row 16: Arabia Saudita
$ws.Range("B16").Value = 298542
$ws.Range("C16").Value = 1227
$ws.Range("D16").Value = 266953
$ws.Range("E16").Value = 28181
$ws.Range("G16").Value = 39
$ws.Range("H16").Value = 3408

# row 24: Irak
$ws.Range("B24").Value = 176931
$ws.Range("C24").Value = 4348
$ws.Range("D24").Value = 125374
$ws.Range("E24").Value = 45697
$ws.Range("G24").Value = 75
$ws.Range("H24").Value = 5860

# row 49: Portugal
$ws.Range("B49").Value = 54102
$ws.Range("C49").Value = 121
$ws.Range("D49").Value = 39697
$ws.Range("E49").Value = 12627
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 1778

# row 61: Uzbekistan
$ws.Range("B61").Value = 35109
$ws.Range("C61").Value = 581
$ws.Range("D61").Value = 30664
$ws.Range("E61").Value = 4215
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = 230

# row 66: Serbia
$ws.Range("B66").Value = 29682
$ws.Range("C66").Value = 211
$ws.Range("D66").Value = 27061
$ws.Range("E66").Value = 1947
$ws.Range("G66").Value = 4
$ws.Range("H66").Value = 674

# row 83: Republica de Macedonia
$ws.Range("B83").Value = 12739
$ws.Range("C83").Value = 86
$ws.Range("D83").Value = 9174
$ws.Range("E83").Value = 3021
$ws.Range("G83").Value = 5
$ws.Range("H83").Value = 544

# row 86: Noruega (D unchanged)
$ws.Range("B86").Value = 9989
$ws.Range("C86").Value = 24
$ws.Range("E86").Value = 871

# row 118: Cuba
$ws.Range("B118").Value = 3316
$ws.Range("C118").Value = 24
$ws.Range("D118").Value = 2620
$ws.Range("E118").Value = 608

# row 174: Islas Feroe (D unchanged)
$ws.Range("B174").Value = 372
$ws.Range("C174").Value = 2
$ws.Range("E174").Value = 147

# rows 213-214: Islas Malvinas and Montserrat swap order/values
# Row 213 now shows "Islas Malvinas" with its updated data
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214 now shows "Montserrat" with its updated data
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
